$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in header values D1:L1 with column letters D..L (shared strings)
$headers = @("D", "E", "F", "G", "H", "I", "J", "K", "L")
$col = 4
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col++
}

# Update the active selection to C1 (matches sheetView selection change)
$ws.Range("C1").Select()
